# MEMBER.xlsx edit: add a new test member row (row 4) to Sheet1,
# matching the "test" / "테스트" / "test@cornsalad.com" / "N" sample data,
# with a mailto hyperlink on the email cell, and move the active selection
# to L4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of sample/test data (mirrors the existing admin/jieun rows above it)
$ws.Range("A4").Value2 = "test"
$ws.Range("B4").Value2 = "test"
$ws.Range("C4").Value2 = "테스트"
$ws.Range("D4").Value2 = "test@cornsalad.com"
[void]$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:test@cornsalad.com")
$ws.Range("I4").Value2 = "N"
$ws.Range("J4").Value2 = "N"

# Move the current selection to L4, as in the saved workbook state
[void]$ws.Range("L4").Select()
